$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting (incl. date number format) from the last existing data row
$ws.Range("A31:C31").Copy()
$ws.Range("A32:C33").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(32, 1).Value = 42341
$ws.Cells.Item(32, 2).Value = 0.5
$ws.Cells.Item(32, 3).Value = "Set up git repository"

$ws.Cells.Item(33, 1).Value = 42355
$ws.Cells.Item(33, 2).Value = 5
$ws.Cells.Item(33, 3).Value = "Make training and test data set to make the classifier"

$ws.Activate()
$ws.Range("B1").Select()
